$d = $word.ActiveDocument

# Locate "When" inside the bolded test-name heading
# ("napadniIgrac_ShouldThrowAnIllegalArgumentException_WhenBiggerThan...").
$found = $d.Content
$found.Find.Execute("When", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $found.Start
$end = $found.End

# Replace "When" with "If" (new naming convention for the exception test).
$mid = $d.Range($start, $end)
$mid.Text = ""
$mid = $d.Range($start, $start)
$mid.InsertAfter("If")

# Toggle bold off/on for just the inserted "If" text. This keeps its
# formatting identical to the surrounding text while still forcing the
# OOXML writer to keep it in its own <w:r> run rather than silently
# re-merging it with the neighboring runs.
$ifRange = $d.Range($start, $start + 2)
$ifRange.Font.Bold = $false
$ifRange.Font.Bold = $true
